# Update of league bases - row data (columns B:AC) got reshuffled between
# matching rows (ids/odds swapped while the row-index column A stayed put).
#
# Pattern observed in the diff:
#   rows 50  <-> 51   (swap)
#   rows 89  <-> 90   (swap)
#   rows 100 -> 103 -> 104 -> 101 -> 102 -> 100  (5-way cycle of data)
#   rows 117 <-> 118  (swap)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture all the "before" values for the affected rows first, so that
# writes to one row never clobber data we still need to read from another. ---
$row50  = $ws.Range("B50:AC50").Value()
$row51  = $ws.Range("B51:AC51").Value()

$row89  = $ws.Range("B89:AC89").Value()
$row90  = $ws.Range("B90:AC90").Value()

$row100 = $ws.Range("B100:AC100").Value()
$row101 = $ws.Range("B101:AC101").Value()
$row102 = $ws.Range("B102:AC102").Value()
$row103 = $ws.Range("B103:AC103").Value()
$row104 = $ws.Range("B104:AC104").Value()

$row117 = $ws.Range("B117:AC117").Value()
$row118 = $ws.Range("B118:AC118").Value()

# --- Simple swaps ---
$ws.Range("B50:AC50").Value   = $row51
$ws.Range("B51:AC51").Value   = $row50

$ws.Range("B89:AC89").Value   = $row90
$ws.Range("B90:AC90").Value   = $row89

$ws.Range("B117:AC117").Value = $row118
$ws.Range("B118:AC118").Value = $row117

# --- 5-way cycle: new(100)=old(102), new(101)=old(104), new(102)=old(101),
#     new(103)=old(100), new(104)=old(103) ---
$ws.Range("B100:AC100").Value = $row102
$ws.Range("B101:AC101").Value = $row104
$ws.Range("B102:AC102").Value = $row101
$ws.Range("B103:AC103").Value = $row100
$ws.Range("B104:AC104").Value = $row103
